{"js": "// Apply the dated-worksheet refresh: update the heading date and each of the\n// 25 two-digit x two-digit multiplication prompts in the table to the next day's set.\nconst replacements = [\n  [\"2025-04-25 Friday\", \"2025-04-26 Saturday\"],\n  [\"65\u00d714=\", \"92\u00d743=\"],\n  [\"83\u00d795=\", \"53\u00d791=\"],\n  [\"69\u00d757=\", \"75\u00d720=\"],\n  [\"67\u00d726=\", \"82\u00d739=\"],\n  [\"72\u00d741=\", \"32\u00d764=\"],\n  [\"60\u00d741=\", \"71\u00d785=\"],\n  [\"33\u00d737=\", \"93\u00d743=\"],\n  [\"92\u00d772=\", \"71\u00d756=\"],\n  [\"71\u00d716=\", \"53\u00d731=\"],\n  [\"65\u00d772=\", \"20\u00d782=\"],\n  [\"73\u00d746=\", \"45\u00d719=\"],\n  [\"62\u00d761=\", \"31\u00d797=\"],\n  [\"85\u00d718=\", \"47\u00d775=\"],\n  [\"94\u00d791=\", \"66\u00d712=\"],\n  [\"49\u00d799=\", \"93\u00d724=\"],\n  [\"91\u00d737=\", \"65\u00d762=\"],\n  [\"73\u00d732=\", \"63\u00d789=\"],\n  [\"46\u00d783=\", \"51\u00d769=\"],\n  [\"79\u00d725=\", \"54\u00d779=\"],\n  [\"68\u00d716=\", \"60\u00d795=\"],\n  [\"32\u00d751=\", \"65\u00d766=\"],\n  [\"18\u00d767=\", \"57\u00d747=\"],\n  [\"25\u00d752=\", \"86\u00d756=\"],\n  [\"84\u00d711=\", \"54\u00d732=\"],\n  [\"30\u00d732=\", \"66\u00d719=\"],\n];\n\nfor (const [find, replace] of replacements) {\n  const results = context.document.body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for: ${find}`);\n  }\n\n  // Each source string is unique in the document, so replace every match found\n  // (normally exactly one run) with the new text, preserving its formatting.\n  for (const range of results.items) {\n    range.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the dated-worksheet refresh: update the heading date and each of the\n# 25 two-digit x two-digit multiplication prompts in the table to the next day's set.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-04-25 Friday\", \"2025-04-26 Saturday\"),\n    @(\"65\u00d714=\", \"92\u00d743=\"),\n    @(\"83\u00d795=\", \"53\u00d791=\"),\n    @(\"69\u00d757=\", \"75\u00d720=\"),\n    @(\"67\u00d726=\", \"82\u00d739=\"),\n    @(\"72\u00d741=\", \"32\u00d764=\"),\n    @(\"60\u00d741=\", \"71\u00d785=\"),\n    @(\"33\u00d737=\", \"93\u00d743=\"),\n    @(\"92\u00d772=\", \"71\u00d756=\"),\n    @(\"71\u00d716=\", \"53\u00d731=\"),\n    @(\"65\u00d772=\", \"20\u00d782=\"),\n    @(\"73\u00d746=\", \"45\u00d719=\"),\n    @(\"62\u00d761=\", \"31\u00d797=\"),\n    @(\"85\u00d718=\", \"47\u00d775=\"),\n    @(\"94\u00d791=\", \"66\u00d712=\"),\n    @(\"49\u00d799=\", \"93\u00d724=\"),\n    @(\"91\u00d737=\", \"65\u00d762=\"),\n    @(\"73\u00d732=\", \"63\u00d789=\"),\n    @(\"46\u00d783=\", \"51\u00d769=\"),\n    @(\"79\u00d725=\", \"54\u00d779=\"),\n    @(\"68\u00d716=\", \"60\u00d795=\"),\n    @(\"32\u00d751=\", \"65\u00d766=\"),\n    @(\"18\u00d767=\", \"57\u00d747=\"),\n    @(\"25\u00d752=\", \"86\u00d756=\"),\n    @(\"84\u00d711=\", \"54\u00d732=\"),\n    @(\"30\u00d732=\", \"66\u00d719=\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    # Fresh Find off the whole document each time so prior replacements don't\n    # shrink/move the search range.\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n    if (-not $found) { throw \"Replace failed for: $oldText\" }\n}\n"}
